$wb = $excel.ActiveWorkbook

# Row 132 (ALC) - hunk 0
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1090.25
$ws.Range("I132").Value = 1143.1
$ws.Range("J132").Value = 297.5
$ws.Range("K132").Value = 3429.3
$ws.Range("L132").Value = 892.5
$ws.Range("M132").Value = -899.2999999999997
$ws.Range("N132").Value = -5952.5

# Row 28 (ARM) - hunk 1
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3235.5
$ws.Range("I28").Value = 3235.5
$ws.Range("K28").Value = 3235.5
$ws.Range("M28").Value = -3043.5

# Row 99 (ARM) - hunk 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 3235.5
$ws.Range("I99").Value = 3235.5
$ws.Range("K99").Value = 3235.5
$ws.Range("M99").Value = -240.5

# Row 132 (ARM) - hunk 3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2394.64
$ws.Range("I132").Value = 2130.318
$ws.Range("K132").Value = 6390.954000000001
$ws.Range("M132").Value = -3860.954000000001

# Row 94 (BSM) - hunk 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1579.8667
$ws.Range("J94").Value = 725
$ws.Range("L94").Value = 725
$ws.Range("N94").Value = -1627

# Row 105 (BSM) - hunk 5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3430.5
$ws.Range("I105").Value = 2495.4348
$ws.Range("K105").Value = 2495.4348
$ws.Range("M105").Value = -748.4348

# Row 11 (CRP) - hunk 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11:N11").ClearContents()

# Row 93 (CRP) - hunk 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 2283.3333
$ws.Range("I93").Value = 2283.3333
$ws.Range("K93").Value = 2283.3333
$ws.Range("M93").Value = -411.3332999999998

# Row 99 (CRP) - hunk 8
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13354.966
$ws.Range("I99").Value = 8480.5
$ws.Range("J99").Value = 16795.766
$ws.Range("K99").Value = 8480.5
$ws.Range("L99").Value = 16795.766
$ws.Range("M99").Value = -6982.5
$ws.Range("N99").Value = -19791.766

# Row 126 (CRP) - hunk 9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 13354.966
$ws.Range("I126").Value = 8480.5
$ws.Range("J126").Value = 16795.766
$ws.Range("K126").Value = 25441.5
$ws.Range("L126").Value = 50387.298
$ws.Range("M126").Value = -22971.5
$ws.Range("N126").Value = -55327.298

# Row 134 (CRP) - hunk 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2595.8845
$ws.Range("I134").Value = 2094.4285
$ws.Range("J134").Value = 4702
$ws.Range("K134").Value = 6283.2855
$ws.Range("L134").Value = 14106
$ws.Range("M134").Value = -3748.2855
$ws.Range("N134").Value = -19176

# Row 113 (CUL) - hunk 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 827.7646999999999
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 833.7143
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 2501.1429
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -6841.1429

# Row 125 (CUL) - hunk 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 129 (CUL) - hunk 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2292.7
$ws.Range("J129").Value = 3672.1667
$ws.Range("L129").Value = 11016.5001
$ws.Range("N129").Value = -21016.5001

# Row 80 (GSM) - hunk 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4358.4736
$ws.Range("I80").Value = 3428.9
$ws.Range("J80").Value = 5391.3335
$ws.Range("K80").Value = 3428.9
$ws.Range("L80").Value = 5391.3335
$ws.Range("M80").Value = -2430.9
$ws.Range("N80").Value = -7387.3335

# Row 83 (GSM) - hunk 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4358.4736
$ws.Range("I83").Value = 3428.9
$ws.Range("J83").Value = 5391.3335
$ws.Range("K83").Value = 17144.5
$ws.Range("L83").Value = 26956.6675
$ws.Range("M83").Value = -12152.5
$ws.Range("N83").Value = -36940.6675

# Row 132 (GSM) - hunk 16
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1463.2222
$ws.Range("I132").Value = 880.7143
$ws.Range("K132").Value = 2642.1429
$ws.Range("M132").Value = -112.1428999999998

# Row 134 (GSM) - hunk 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 95851.86
$ws.Range("J134").Value = 95851.86
$ws.Range("L134").Value = 287555.58
$ws.Range("N134").Value = -292625.58

# Row 2 (LTW) - hunk 18
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20000000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 7 (LTW) - hunk 19
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1489.7646
$ws.Range("I7").Value = 1609.3
$ws.Range("J7").Value = 1319
$ws.Range("K7").Value = 1609.3
$ws.Range("L7").Value = 1319
$ws.Range("M7").Value = -1497.3
$ws.Range("N7").Value = -1543

# Row 22 (LTW) - hunk 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5222
$ws.Range("I22").Value = 1528.25
$ws.Range("K22").Value = 1528.25
$ws.Range("M22").Value = -1233.25

# Row 27 (LTW) - hunk 21
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5222
$ws.Range("I27").Value = 1528.25
$ws.Range("K27").Value = 1528.25
$ws.Range("M27").Value = -1421.25

# Row 61 (LTW) - hunk 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2607.261
$ws.Range("I61").Value = 3024.1333
$ws.Range("J61").Value = 1825.625
$ws.Range("K61").Value = 3024.1333
$ws.Range("L61").Value = 1825.625
$ws.Range("M61").Value = -2822.1333
$ws.Range("N61").Value = -2229.625

# Row 82 (LTW) - hunk 23
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5101
$ws.Range("I82").Value = 5603.6665
$ws.Range("J82").Value = 4849.6665
$ws.Range("K82").Value = 5603.6665
$ws.Range("L82").Value = 4849.6665
$ws.Range("M82").Value = -5242.6665
$ws.Range("N82").Value = -5571.6665

# Row 85 (LTW) - hunk 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5101
$ws.Range("I85").Value = 5603.6665
$ws.Range("J85").Value = 4849.6665
$ws.Range("K85").Value = 5603.6665
$ws.Range("L85").Value = 4849.6665
$ws.Range("M85").Value = -4355.6665
$ws.Range("N85").Value = -7345.6665

# Row 100 (LTW) - hunk 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5538.8
$ws.Range("I100").Value = 673
$ws.Range("J100").Value = 12837.5
$ws.Range("K100").Value = 673
$ws.Range("L100").Value = 12837.5
$ws.Range("M100").Value = -132
$ws.Range("N100").Value = -13919.5

# Row 113 (LTW) - hunk 26
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2607.261
$ws.Range("I113").Value = 3024.1333
$ws.Range("J113").Value = 1825.625
$ws.Range("K113").Value = 3024.1333
$ws.Range("L113").Value = 1825.625
$ws.Range("M113").Value = -854.1333
$ws.Range("N113").Value = -6165.625

# Row 126 (LTW) - hunk 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1489.7646
$ws.Range("I126").Value = 1609.3
$ws.Range("J126").Value = 1319
$ws.Range("K126").Value = 4827.9
$ws.Range("L126").Value = 3957
$ws.Range("M126").Value = -2357.9
$ws.Range("N126").Value = -8897

# Row 132 (LTW) - hunk 28
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4860.5
$ws.Range("I132").Value = 4739
$ws.Range("J132").Value = 4928
$ws.Range("K132").Value = 14217
$ws.Range("L132").Value = 14784
$ws.Range("M132").Value = -11687
$ws.Range("N132").Value = -19844

# Row 136 (LTW) - hunk 29
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1914.5
$ws.Range("J136").Value = 1488.6666
$ws.Range("L136").Value = 4465.9998
$ws.Range("N136").Value = -9565.9998

# Row 107 (WVR) - hunk 30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 560.8125
$ws.Range("I107").Value = 363.125
$ws.Range("J107").Value = 758.5
$ws.Range("K107").Value = 1089.375
$ws.Range("L107").Value = 2275.5
$ws.Range("M107").Value = 830.625
$ws.Range("N107").Value = -6115.5

# Row 132 (WVR) - hunk 31
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132:N132").ClearContents()
